$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values of rows 2, 4, 5 (row 3 is unchanged)
$row2 = @{
    D = $ws.Range("D2").Value2
    J = $ws.Range("J2").Value2
    K = $ws.Range("K2").Value2
    L = $ws.Range("L2").Value2
    M = $ws.Range("M2").Value2
    P = $ws.Range("P2").Value2
}
$row4 = @{
    D = $ws.Range("D4").Value2
    J = $ws.Range("J4").Value2
    K = $ws.Range("K4").Value2
    L = $ws.Range("L4").Value2
    M = $ws.Range("M4").Value2
    P = $ws.Range("P4").Value2
}
$row5 = @{
    D = $ws.Range("D5").Value2
    J = $ws.Range("J5").Value2
    K = $ws.Range("K5").Value2
    L = $ws.Range("L5").Value2
    M = $ws.Range("M5").Value2
    P = $ws.Range("P5").Value2
}

# Row 2 gets the old values of row 4
$ws.Range("D2").Value2 = $row4.D
$ws.Range("J2").Value2 = $row4.J
$ws.Range("K2").Value2 = $row4.K
$ws.Range("L2").Value2 = $row4.L
$ws.Range("M2").Value2 = $row4.M
$ws.Range("P2").Value2 = $row4.P

# Row 4 gets the old values of row 5
$ws.Range("D4").Value2 = $row5.D
$ws.Range("J4").Value2 = $row5.J
$ws.Range("K4").Value2 = $row5.K
$ws.Range("L4").Value2 = $row5.L
$ws.Range("M4").Value2 = $row5.M
$ws.Range("P4").Value2 = $row5.P

# Row 5 gets the old values of row 2
$ws.Range("D5").Value2 = $row2.D
$ws.Range("J5").Value2 = $row2.J
$ws.Range("K5").Value2 = $row2.K
$ws.Range("L5").Value2 = $row2.L
$ws.Range("M5").Value2 = $row2.M
$ws.Range("P5").Value2 = $row2.P

$wb.Save()
